# Apply the localization-sheet edits:
#  - bump the "Quest" id values from 20000 to 20001 on both quest sheets
#  - widen column A on the QuestNameEntities sheet
#  - update the saved cursor/selection + active-tab state to match what was
#    left selected when the file was last saved (NameEntities tab active,
#    with C15 selected; QuestNameEntities/QuestDescriptionEntities left with
#    A2 selected)

$wb = $excel.ActiveWorkbook

# --- QuestNameEntities ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = 20001
$ws3.Columns.Item(1).ColumnWidth = 20.25
$ws3.Range("A2").Select() | Out-Null

# --- QuestDescriptionEntities ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = 20001
$ws4.Range("A2").Select() | Out-Null

# --- NameEntities becomes the active sheet/selection last, matching the diff ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C15").Select() | Out-Null
